# Apply "track replacement schedule" updates.
# - "Operating Costs" sheet: column F ("Tire/track ($/hr)") values updated
#   to reflect the new track replacement schedule.
# - "Summary Costs" sheet: column I ("Use Cost ($/PMH)") values recomputed
#   to reflect the new Tire/track costs.

$wb = $excel.ActiveWorkbook

# --- Operating Costs sheet: update Tire/track ($/hr) column (F) ---
$wsOps = $wb.Worksheets.Item("Operating Costs")

$wsOps.Range("F2").Value = 4.025
$wsOps.Range("F4").Value = 2.3
$wsOps.Range("F5").Value = 2.3
$wsOps.Range("F6").Value = 2.68333333333
$wsOps.Range("F7").Value = 2.3
$wsOps.Range("F8").Value = 4.3125
$wsOps.Range("F9").Value = 2.15625

# --- Summary Costs sheet: update Use Cost ($/PMH) column (I) ---
$wsSummary = $wb.Worksheets.Item("Summary Costs")

$wsSummary.Range("I2").Value = 74.3185879997
$wsSummary.Range("I4").Value = 76.4702978886
$wsSummary.Range("I5").Value = 44.9658100045
$wsSummary.Range("I6").Value = 59.3833116995
$wsSummary.Range("I7").Value = 44.5127603491
$wsSummary.Range("I8").Value = 60.505122612
$wsSummary.Range("I9").Value = 101.094962625
